$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes: A2=0 (unchanged), B2="" (text, empty), C2=<hash> (text)
$ws.Range("A2").Value = 0

# A plain "" assignment clears the cell entirely (drops its text type).
# Forcing a leading apostrophe keeps the cell typed as text while being
# empty, then resetting the style clears the quote-prefix formatting that
# the apostrophe entry would otherwise leave behind.
$ws.Range("B2").Value = "'"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "5237938ceb8907ff563b787544e65339"
